$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "https://www.autodoc.co.uk/ridex/8095160"
$ws.Range("D2").Value = "£7. 59"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-42%"

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "20%"
